$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 26
$ws.Range("H3").Value = 2
$ws.Range("G4").Value = 2

$ws.Range("F4").Select()
